# Sample Project "Rules" sheet: row 11 (R40 rule) Greeting-to value (B11)
# changes from the text "R40" to the text "1". We write it through a
# TEXT() formula and then paste-special just the value back over itself
# so the result lands as a literal (non-formula) string, keeping the
# cell's existing style/border formatting untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$cell = $ws.Range("B11")
$cell.Formula = '=TEXT(1,"0")'
$cell.Copy()
$cell.PasteSpecial(-4163)
